# Add a new "2022-Q4" sheet (fund-holdings detail) right before the
# existing "2022-Q3" sheet, and update the "总计" (summary) sheet with a
# new leading row for 2022-Q4, pushing every other quarter down by one.

$wb = $excel.ActiveWorkbook

function Set-IndexStyle($cell) {
    # Replicates the bold / thin-boxed / centered "index column" look
    # used by the other index (A) and header cells in this workbook.
    $cell.Font.Bold = $true
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet before "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($q3)
$newSheet.Name = "2022-Q4"

# Header row (B1:H1)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $cell = $newSheet.Cells.Item(1, $col)
    $cell.Value = $h
    Set-IndexStyle $cell
    $col = $col + 1
}

# Columns B (fund code) and D:G (scale/position/value) are stored as
# text in the source workbook (they keep leading zeros / trailing
# zeros), so force a text number format before writing them.
$newSheet.Range("B2:B12").NumberFormat = "@"
$newSheet.Range("D2:G12").NumberFormat = "@"

$rows = @(
    @(0,  "001040", "新华策略精选股票",                 "13.52", "94.19", "3.38", "0.4570", 8),
    @(1,  "519087", "新华优选分红混合",                 "11.44", "88.79", "3.95", "0.4519", 7),
    @(2,  "519156", "新华行业轮换灵活配置混合A",         "6.37",  "94.38", "3.82", "0.2433", 7),
    @(3,  "519158", "新华趋势领航混合",                 "5.42",  "94.63", "3.17", "0.1718", 9),
    @(4,  "004982", "新华安享多裕定期开放灵活配置混合",   "3.13",  "45.81", "2.93", "0.0917", 4),
    @(5,  "001294", "新华战略新兴产业灵活配置混合",       "1.24",  "94.39", "5.04", "0.0625", 4),
    @(6,  "011457", "新华行业龙头主题股票",             "0.94",  "94.57", "4.56", "0.0429", 7),
    @(7,  "005293", "诺德新旺灵活配置混合",             "0.55",  "92.81", "6.22", "0.0342", 3),
    @(8,  "519157", "新华行业轮换灵活配置混合C",         "0.83",  "94.38", "3.82", "0.0317", 7),
    @(9,  "005209", "东吴双三角股票A",                 "0.09",  "92.11", "3.11", "0.0028", 9),
    @(10, "005210", "东吴双三角股票C",                 "0.09",  "92.11", "3.11", "0.0028", 9)
)

$r = 2
foreach ($row in $rows) {
    $aCell = $newSheet.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    Set-IndexStyle $aCell

    $newSheet.Range("B$r").Value = $row[1]
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("D$r").Value = $row[3]
    $newSheet.Range("E$r").Value = $row[4]
    $newSheet.Range("F$r").Value = $row[5]
    $newSheet.Range("G$r").Value = $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: add the 2022-Q4 row at the top of
#    the data block and shift the existing quarters down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @(0, "2022-Q4", 11, 1.59),
    @(1, "2022-Q3", 15, 2.51),
    @(2, "2022-Q2", 7,  0.89),
    @(3, "2022-Q1", 3,  0.25),
    @(4, "2021-Q4", 1,  0.16),
    @(5, "2021-Q2", 15, 0.58),
    @(6, "2021-Q1", 3,  0.03),
    @(7, "2020-Q4", 3,  0.04)
)

$r = 2
foreach ($row in $summaryRows) {
    $aCell = $summary.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    if ($r -eq 9) {
        # Row 9 is brand new (sheet previously ended at row 8) so it
        # won't inherit the existing index-column formatting.
        Set-IndexStyle $aCell
    }

    $summary.Range("B$r").Value = $row[1]
    $summary.Range("C$r").Value = $row[2]
    $summary.Range("D$r").Value = $row[3]
    $r = $r + 1
}

Write-Output "done"
